{"js": "// Replace the shorthand inequality symbols used for fish-length cutoffs\n// with their spelled-out equivalents:\n//   \"fry (<46 mm FL)\"                 -> \"fry (less than 46 mm FL)\"\n//   \"pre-smolt/smolt passage (>45 mm TL)\" -> \"pre-smolt/smolt passage (greater than 45 mm TL)\"\nconst body = context.document.body;\n\nconst oldText =\n  \"fish passage as total passage along with fry (<46 mm FL) and pre-smolt/smolt passage (>45 mm TL) for all four runs of Chinook\";\nconst newText =\n  \"fish passage as total passage along with fry (less than 46 mm FL) and pre-smolt/smolt passage (greater than 45 mm TL) for all four runs of Chinook\";\n\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n} else {\n  // Fallback: perform the two substitutions independently in case the\n  // longer phrase above no longer matches verbatim (e.g. due to other\n  // nearby edits), still targeting only the intended \"<\"/\">\" occurrences.\n  const ltResults = body.search(\"(<46 mm FL)\", { matchCase: true });\n  ltResults.load(\"items\");\n  await context.sync();\n  if (ltResults.items.length > 0) {\n    ltResults.items[0].insertText(\"(less than 46 mm FL)\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  const gtResults = body.search(\"(>45 mm TL)\", { matchCase: true });\n  gtResults.load(\"items\");\n  await context.sync();\n  if (gtResults.items.length > 0) {\n    gtResults.items[0].insertText(\"(greater than 45 mm TL)\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Replace the shorthand inequality symbols used for fish-length cutoffs\n# with their spelled-out equivalents:\n#   \"fry (<46 mm FL)\"                     -> \"fry (less than 46 mm FL)\"\n#   \"pre-smolt/smolt passage (>45 mm TL)\" -> \"pre-smolt/smolt passage (greater than 45 mm TL)\"\n\n$d = $word.ActiveDocument\n\n$oldText = \"fish passage as total passage along with fry (<46 mm FL) and pre-smolt/smolt passage (>45 mm TL) for all four runs of Chinook\"\n$newText = \"fish passage as total passage along with fry (less than 46 mm FL) and pre-smolt/smolt passage (greater than 45 mm TL) for all four runs of Chinook\"\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Replacement.ClearFormatting()\n$found = $range.Find.Execute(\n    $oldText,\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    $newText,\n    2\n)\n\nif (-not $found) {\n    # Fallback: perform the two substitutions independently in case the\n    # longer phrase above no longer matches verbatim (e.g. due to other\n    # nearby edits), still targeting only the intended \"<\" / \">\" spots.\n    $ltRange = $d.Content\n    $ltRange.Find.ClearFormatting()\n    $ltRange.Find.Replacement.ClearFormatting()\n    $ltRange.Find.Execute(\n        \"(<46 mm FL)\",\n        $false, $false, $false, $false, $false, $true, 1, $false,\n        \"(less than 46 mm FL)\",\n        2\n    ) | Out-Null\n\n    $gtRange = $d.Content\n    $gtRange.Find.ClearFormatting()\n    $gtRange.Find.Replacement.ClearFormatting()\n    $gtRange.Find.Execute(\n        \"(>45 mm TL)\",\n        $false, $false, $false, $false, $false, $true, 1, $false,\n        \"(greater than 45 mm TL)\",\n        2\n    ) | Out-Null\n}\n"}
